$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns A and D (date-like / numeric-looking IDs) are text before assignment
$ws.Range("A2:A22").NumberFormat = "@"
$ws.Range("D2:D22").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = '2025-07-16'
$ws.Cells.Item(2, 2).Value = 60
$ws.Cells.Item(2, 3).Value = 'V V REFEICOES LTDA'
$ws.Cells.Item(2, 4).Value = '54396269'
$ws.Cells.Item(2, 5).Value = 42173656
$ws.Cells.Item(2, 6).Value = 'COADOR DE CAFÉ G'
$ws.Cells.Item(2, 7).Value = 23
$ws.Cells.Item(2, 8).Value = 6.07
$ws.Cells.Item(2, 9).Value = 15.64

$ws.Cells.Item(3, 1).Value = '2025-07-16'
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = 'V V REFEICOES LTDA'
$ws.Cells.Item(3, 4).Value = '54396269'
$ws.Cells.Item(3, 5).Value = 17541022
$ws.Cells.Item(3, 6).Value = 'COADOR DE CAFE INDUSTRIAL (MAIOR)'
$ws.Cells.Item(3, 7).Value = 26
$ws.Cells.Item(3, 8).Value = 2.12
$ws.Cells.Item(3, 9).Value = 2.12

$ws.Cells.Item(4, 1).Value = '2025-07-16'
$ws.Cells.Item(4, 2).Value = 10
$ws.Cells.Item(4, 3).Value = 'V V REFEICOES LTDA'
$ws.Cells.Item(4, 4).Value = '54396269'
$ws.Cells.Item(4, 5).Value = 17125814
$ws.Cells.Item(4, 6).Value = 'COADOR DE CAFE P'
$ws.Cells.Item(4, 7).Value = 21
$ws.Cells.Item(4, 8).Value = 2.18
$ws.Cells.Item(4, 9).Value = 2.04

$ws.Cells.Item(5, 1).Value = '2025-07-16'
$ws.Cells.Item(5, 2).Value = 70
$ws.Cells.Item(5, 3).Value = 'TECHLOG - SERVICOS DE GESTAO E SISTEMAS INFORMATIZ'
$ws.Cells.Item(5, 4).Value = '54316914'
$ws.Cells.Item(5, 5).Value = 19264853
$ws.Cells.Item(5, 6).Value = 'PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE'
$ws.Cells.Item(5, 7).Value = 722
$ws.Cells.Item(5, 8).Value = 8.210000000000001
$ws.Cells.Item(5, 9).Value = 13.48

$ws.Cells.Item(6, 1).Value = '2025-07-16'
$ws.Cells.Item(6, 2).Value = 40
$ws.Cells.Item(6, 3).Value = 'CONDOMINIO RESIDENCIAL EPHYGENIO SALLES'
$ws.Cells.Item(6, 4).Value = '54365832'
$ws.Cells.Item(6, 5).Value = 20619556
$ws.Cells.Item(6, 6).Value = 'PEDRA SANITARIA 35G RUBI FLORAL'
$ws.Cells.Item(6, 7).Value = 118
$ws.Cells.Item(6, 8).Value = 13.21
$ws.Cells.Item(6, 9).Value = 8.550000000000001

$ws.Cells.Item(7, 1).Value = '2025-07-16'
$ws.Cells.Item(7, 2).Value = 10
$ws.Cells.Item(7, 3).Value = 'CONDOMINIO DO TVLANDIA MALL'
$ws.Cells.Item(7, 4).Value = '54119372'
$ws.Cells.Item(7, 5).Value = 28133466
$ws.Cells.Item(7, 6).Value = 'DESINFETANTE CONCENTRADO 5L AUDAX MAX 1:200 - LAVANDA'
$ws.Cells.Item(7, 7).Value = 105
$ws.Cells.Item(7, 8).Value = 2.31
$ws.Cells.Item(7, 9).Value = 1.83

$ws.Cells.Item(8, 1).Value = '2025-07-17'
$ws.Cells.Item(8, 2).Value = 100
$ws.Cells.Item(8, 3).Value = 'TEL TELECOMUNICACOES LTDA.'
$ws.Cells.Item(8, 4).Value = '54443314'
$ws.Cells.Item(8, 5).Value = 11936640
$ws.Cells.Item(8, 6).Value = 'LIMPADOR VEJA MULTIUSO GOLD 500ML'
$ws.Cells.Item(8, 7).Value = 3312
$ws.Cells.Item(8, 8).Value = 10.39
$ws.Cells.Item(8, 9).Value = 14.42

$ws.Cells.Item(9, 1).Value = '2025-07-21'
$ws.Cells.Item(9, 2).Value = 14832
$ws.Cells.Item(9, 3).Value = 'BMS INDUSTRIA E COMERCIO DE PRODUTOS ALIMENTICIOS E DISTRIBU'
$ws.Cells.Item(9, 4).Value = '54586025'
$ws.Cells.Item(9, 5).Value = 11939543
$ws.Cells.Item(9, 6).Value = 'RODO COM CABO M 40CM'
$ws.Cells.Item(9, 7).Value = 7
$ws.Cells.Item(9, 8).Value = 44.98
$ws.Cells.Item(9, 9).Value = 780.46

$ws.Cells.Item(10, 1).Value = '2025-07-21'
$ws.Cells.Item(10, 2).Value = 35000
$ws.Cells.Item(10, 3).Value = 'BMS INDUSTRIA E COMERCIO DE PRODUTOS ALIMENTICIOS E DISTRIBU'
$ws.Cells.Item(10, 4).Value = '54586025'
$ws.Cells.Item(10, 5).Value = 12067332
$ws.Cells.Item(10, 6).Value = 'PANO DE CHAO BRANCO ALVEJADO CC COMUM 43X70CM'
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 2830.62
$ws.Cells.Item(10, 9).Value = 9667.6

$ws.Cells.Item(11, 1).Value = '2025-07-21'
$ws.Cells.Item(11, 2).Value = 2400
$ws.Cells.Item(11, 3).Value = 'BMS INDUSTRIA E COMERCIO DE PRODUTOS ALIMENTICIOS E DISTRIBU'
$ws.Cells.Item(11, 4).Value = '54586025'
$ws.Cells.Item(11, 5).Value = 11939672
$ws.Cells.Item(11, 6).Value = 'VASSOURA PIACAVA 20 FUROS'
$ws.Cells.Item(11, 7).Value = -2131
$ws.Cells.Item(11, 8).Value = 16.57
$ws.Cells.Item(11, 9).Value = 128.26

$ws.Cells.Item(12, 1).Value = '2025-07-21'
$ws.Cells.Item(12, 2).Value = 200
$ws.Cells.Item(12, 3).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(12, 4).Value = '54503121'
$ws.Cells.Item(12, 5).Value = 13996941
$ws.Cells.Item(12, 6).Value = 'SACO DE LIXO 30L REFORCADO PACOTINHO C/10 UND FORTE MAX'
$ws.Cells.Item(12, 7).Value = 81
$ws.Cells.Item(12, 8).Value = 13.52
$ws.Cells.Item(12, 9).Value = 22.52

$ws.Cells.Item(13, 1).Value = '2025-07-22'
$ws.Cells.Item(13, 2).Value = 22
$ws.Cells.Item(13, 3).Value = 'CONDOMINIO CRISTAL TOWER'
$ws.Cells.Item(13, 4).Value = '54625981'
$ws.Cells.Item(13, 5).Value = 11939645
$ws.Cells.Item(13, 6).Value = 'SACO DE LIXO 50L PRETO COMUM - PCT C/100 UND'
$ws.Cells.Item(13, 7).Value = 153
$ws.Cells.Item(13, 8).Value = 4.3
$ws.Cells.Item(13, 9).Value = 4.81

$ws.Cells.Item(14, 1).Value = '2025-07-22'
$ws.Cells.Item(14, 2).Value = 9
$ws.Cells.Item(14, 3).Value = 'CARITAS ARQUIDIOCESANA DE MANAUS'
$ws.Cells.Item(14, 4).Value = '54625975'
$ws.Cells.Item(14, 5).Value = 17059594
$ws.Cells.Item(14, 6).Value = 'SABAO EM PO ESPUMIL 4KG'
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 1.91
$ws.Cells.Item(14, 9).Value = 1.85

$ws.Cells.Item(15, 1).Value = '2025-07-23'
$ws.Cells.Item(15, 2).Value = 400
$ws.Cells.Item(15, 3).Value = 'V V REFEICOES LTDA'
$ws.Cells.Item(15, 4).Value = '54692772'
$ws.Cells.Item(15, 5).Value = 32130390
$ws.Cells.Item(15, 6).Value = 'ESPONJA MULTIUSO JEITOSA'
$ws.Cells.Item(15, 7).Value = 11606
$ws.Cells.Item(15, 8).Value = 21.42
$ws.Cells.Item(15, 9).Value = 51.23

$ws.Cells.Item(16, 1).Value = '2025-07-24'
$ws.Cells.Item(16, 2).Value = 40
$ws.Cells.Item(16, 3).Value = 'MUNDY SERVIS LIMPEZAS PREDIAIS LTDA'
$ws.Cells.Item(16, 4).Value = '54704921'
$ws.Cells.Item(16, 5).Value = 20619556
$ws.Cells.Item(16, 6).Value = 'PEDRA SANITARIA 35G RUBI FLORAL'
$ws.Cells.Item(16, 7).Value = 118
$ws.Cells.Item(16, 8).Value = 13.21
$ws.Cells.Item(16, 9).Value = 8.550000000000001

$ws.Cells.Item(17, 1).Value = '2025-07-24'
$ws.Cells.Item(17, 2).Value = 70
$ws.Cells.Item(17, 3).Value = 'SUPER G - TRANSPORTE E LOCACAO DE EQUIPAMENTOS INDUSTRIAIS L'
$ws.Cells.Item(17, 4).Value = '54736983'
$ws.Cells.Item(17, 5).Value = 12097805
$ws.Cells.Item(17, 6).Value = 'FLANELA LARANJA TAM G 38x58CM'
$ws.Cells.Item(17, 7).Value = -11
$ws.Cells.Item(17, 8).Value = 15.7
$ws.Cells.Item(17, 9).Value = 13.96

$ws.Cells.Item(18, 1).Value = '2025-07-25'
$ws.Cells.Item(18, 2).Value = 100
$ws.Cells.Item(18, 3).Value = 'MANJAR SERVICOS GERAIS SA'
$ws.Cells.Item(18, 4).Value = '54799813'
$ws.Cells.Item(18, 5).Value = 31186309
$ws.Cells.Item(18, 6).Value = 'LA DE ACO ASSOLAN 45G'
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).Value = 9.44
$ws.Cells.Item(18, 9).Value = 17.62

$ws.Cells.Item(19, 1).Value = '2025-07-28'
$ws.Cells.Item(19, 2).Value = 30
$ws.Cells.Item(19, 3).Value = 'PARENTE ANDRADE LTDA'
$ws.Cells.Item(19, 4).Value = '54869922'
$ws.Cells.Item(19, 5).Value = 14413867
$ws.Cells.Item(19, 6).Value = 'AGUA SANITARIA 5L GLOBO SAN'
$ws.Cells.Item(19, 7).Value = 63
$ws.Cells.Item(19, 8).Value = 4.31
$ws.Cells.Item(19, 9).Value = 8.41

$ws.Cells.Item(20, 1).Value = '2025-07-29'
$ws.Cells.Item(20, 2).Value = 85
$ws.Cells.Item(20, 3).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(20, 4).Value = '54402873'
$ws.Cells.Item(20, 5).Value = 11938367
$ws.Cells.Item(20, 6).Value = 'PAPEL HIGIENICO 8X300 NEWPAPER 100% Celulose'
$ws.Cells.Item(20, 7).Value = 293
$ws.Cells.Item(20, 8).Value = 4.89
$ws.Cells.Item(20, 9).Value = 6.86

$ws.Cells.Item(21, 1).Value = '2025-07-29'
$ws.Cells.Item(21, 2).Value = 168
$ws.Cells.Item(21, 3).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(21, 4).Value = '54402873'
$ws.Cells.Item(21, 5).Value = 15011531
$ws.Cells.Item(21, 6).Value = 'DETERGENTE LIMPOL COCO 500ML'
$ws.Cells.Item(21, 7).Value = 133
$ws.Cells.Item(21, 8).Value = 21.89
$ws.Cells.Item(21, 9).Value = 34.09

$ws.Cells.Item(22, 1).Value = '2025-07-29'
$ws.Cells.Item(22, 2).Value = 350
$ws.Cells.Item(22, 3).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(22, 4).Value = '54402873'
$ws.Cells.Item(22, 5).Value = 12054191
$ws.Cells.Item(22, 6).Value = 'PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM'
$ws.Cells.Item(22, 7).Value = 4835
$ws.Cells.Item(22, 8).Value = 20.15
$ws.Cells.Item(22, 9).Value = 39.65

# Reset style on A and D columns back to Normal to drop temporary text formatting marker
$ws.Range("A2:A22").Style = "Normal"
$ws.Range("D2:D22").Style = "Normal"

